# Update stats for 2025-09
# Append a new row (row 22) to Sheet1 with the September 2025 statistics,
# matching the date-formatted style used in column A of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 22

$ws.Range("A$newRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A$newRow").Value = 45901
$ws.Range("B$newRow").Value = 6269
$ws.Range("C$newRow").Value = 992
$ws.Range("D$newRow").Value = 5696013
$ws.Range("E$newRow").Value = 908.5999361939703
$ws.Range("F$newRow").Value = 7.918746772249952
$ws.Range("G$newRow").Value = 3.765690376569042
$ws.Range("H$newRow").Value = 23.87053210274967
